$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the filename shared string (used by B1 and A2)
$ws.Range("B1").Value = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\net2.pdf"
$ws.Range("A2").Value = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\net2.pdf"

# Update similarity score values
$ws.Range("B2").Value = 0.9999999403953509
$ws.Range("C2").Value = 0.4796421785405556
$ws.Range("B3").Value = 0.4796421785405556
$ws.Range("C3").Value = 0.9999999999999998
